# Apply scheduled market-data refresh updates to Leve profit tables
# across all profession sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1487.5428
$ws.Cells.Item(112, 9).Value = 661.25
$ws.Cells.Item(112, 10).Value = 1732.3704
$ws.Cells.Item(112, 11).Value = 1983.75
$ws.Cells.Item(112, 12).Value = 5197.1112
$ws.Cells.Item(112, 13).Value = -875.75
$ws.Cells.Item(112, 14).Value = -7413.1112
$ws.Cells.Item(132, 8).Value = 989.617
$ws.Cells.Item(132, 9).Value = 815.425
$ws.Cells.Item(132, 10).Value = 1985
$ws.Cells.Item(132, 11).Value = 2446.275
$ws.Cells.Item(132, 12).Value = 5955
$ws.Cells.Item(132, 13).Value = 83.72500000000036
$ws.Cells.Item(132, 14).Value = -11015
$ws.Cells.Item(135, 8).Value = 2325.8865
$ws.Cells.Item(135, 9).Value = 1067
$ws.Cells.Item(135, 11).Value = 9603
$ws.Cells.Item(135, 13).Value = -7068
$ws.Cells.Item(137, 8).Value = 873.6977000000001
$ws.Cells.Item(137, 9).Value = 763.9545000000001
$ws.Cells.Item(137, 10).Value = 988.6667
$ws.Cells.Item(137, 11).Value = 2291.8635
$ws.Cells.Item(137, 12).Value = 2966.0001
$ws.Cells.Item(137, 13).Value = 258.1364999999996
$ws.Cells.Item(137, 14).Value = -8066.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 356216.84
$ws.Cells.Item(32, 9).Value = 2756.808
$ws.Cells.Item(32, 11).Value = 2756.808
$ws.Cells.Item(32, 13).Value = -2469.808
$ws.Cells.Item(61, 8).Value = 1322.2572
$ws.Cells.Item(61, 9).Value = 924.0476
$ws.Cells.Item(61, 10).Value = 1919.5714
$ws.Cells.Item(61, 11).Value = 924.0476
$ws.Cells.Item(61, 12).Value = 1919.5714
$ws.Cells.Item(61, 13).Value = -712.0476
$ws.Cells.Item(61, 14).Value = -2343.5714
$ws.Cells.Item(74, 8).Value = 1326.6571
$ws.Cells.Item(74, 9).Value = 1312.6296
$ws.Cells.Item(74, 10).Value = 1374
$ws.Cells.Item(74, 11).Value = 1312.6296
$ws.Cells.Item(74, 12).Value = 1374
$ws.Cells.Item(74, 13).Value = -438.6296
$ws.Cells.Item(74, 14).Value = -3122
$ws.Cells.Item(77, 8).Value = 1326.6571
$ws.Cells.Item(77, 9).Value = 1312.6296
$ws.Cells.Item(77, 10).Value = 1374
$ws.Cells.Item(77, 11).Value = 6563.148
$ws.Cells.Item(77, 12).Value = 6870
$ws.Cells.Item(77, 13).Value = -2195.148
$ws.Cells.Item(77, 14).Value = -15606
$ws.Cells.Item(97, 8).Value = 1383.3226
$ws.Cells.Item(97, 9).Value = 1069.5294
$ws.Cells.Item(97, 10).Value = 1764.3572
$ws.Cells.Item(97, 11).Value = 1069.5294
$ws.Cells.Item(97, 12).Value = 1764.3572
$ws.Cells.Item(97, 13).Value = -573.5293999999999
$ws.Cells.Item(97, 14).Value = -2756.3572
$ws.Cells.Item(136, 8).Value = 1322.2572
$ws.Cells.Item(136, 9).Value = 924.0476
$ws.Cells.Item(136, 10).Value = 1919.5714
$ws.Cells.Item(136, 11).Value = 2772.1428
$ws.Cells.Item(136, 12).Value = 5758.7142
$ws.Cells.Item(136, 13).Value = -222.1428000000001
$ws.Cells.Item(136, 14).Value = -10858.7142

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 5893.946
$ws.Cells.Item(134, 9).Value = 1649.3334
$ws.Cells.Item(134, 10).Value = 13730.154
$ws.Cells.Item(134, 11).Value = 4948.0002
$ws.Cells.Item(134, 12).Value = 41190.462
$ws.Cells.Item(134, 13).Value = -2413.0002
$ws.Cells.Item(134, 14).Value = -46260.462

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7939453
$ws.Cells.Item(31, 9).Value = 9805528
$ws.Cells.Item(31, 10).Value = 8635.333000000001
$ws.Cells.Item(31, 11).Value = 9805528
$ws.Cells.Item(31, 12).Value = 8635.333000000001
$ws.Cells.Item(31, 13).Value = -9805233
$ws.Cells.Item(31, 14).Value = -9225.333000000001
$ws.Cells.Item(34, 8).Value = 7939453
$ws.Cells.Item(34, 9).Value = 9805528
$ws.Cells.Item(34, 10).Value = 8635.333000000001
$ws.Cells.Item(34, 11).Value = 9805528
$ws.Cells.Item(34, 12).Value = 8635.333000000001
$ws.Cells.Item(34, 13).Value = -9805326
$ws.Cells.Item(34, 14).Value = -9039.333000000001
$ws.Cells.Item(58, 8).Value = 1432.9678
$ws.Cells.Item(58, 9).Value = 1211.4546
$ws.Cells.Item(58, 10).Value = 1974.4445
$ws.Cells.Item(58, 11).Value = 1211.4546
$ws.Cells.Item(58, 12).Value = 1974.4445
$ws.Cells.Item(58, 13).Value = -1008.4546
$ws.Cells.Item(58, 14).Value = -2380.4445
$ws.Cells.Item(132, 8).Value = 1604.1282
$ws.Cells.Item(132, 9).Value = 1361.0571
$ws.Cells.Item(132, 10).Value = 3731
$ws.Cells.Item(132, 11).Value = 4083.1713
$ws.Cells.Item(132, 12).Value = 11193
$ws.Cells.Item(132, 13).Value = -1553.1713
$ws.Cells.Item(132, 14).Value = -16253
$ws.Cells.Item(134, 8).Value = 1088.629
$ws.Cells.Item(134, 9).Value = 1033.8572
$ws.Cells.Item(134, 11).Value = 3101.5716
$ws.Cells.Item(134, 13).Value = -566.5715999999998
$ws.Cells.Item(136, 8).Value = 1432.9678
$ws.Cells.Item(136, 9).Value = 1211.4546
$ws.Cells.Item(136, 10).Value = 1974.4445
$ws.Cells.Item(136, 11).Value = 3634.3638
$ws.Cells.Item(136, 12).Value = 5923.333500000001
$ws.Cells.Item(136, 13).Value = -1084.3638
$ws.Cells.Item(136, 14).Value = -11023.3335

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 114.47619
$ws.Cells.Item(12, 10).Value = 120.15
$ws.Cells.Item(12, 12).Value = 360.45
$ws.Cells.Item(12, 14).Value = -706.45
$ws.Cells.Item(69, 8).Value = 1650
$ws.Cells.Item(72, 8).Value = 1650
$ws.Cells.Item(129, 8).Value = 15152734
$ws.Cells.Item(129, 9).Value = 1246
$ws.Cells.Item(129, 10).Value = 19609054
$ws.Cells.Item(129, 11).Value = 3738
$ws.Cells.Item(129, 12).Value = 58827162
$ws.Cells.Item(129, 13).Value = 1262
$ws.Cells.Item(129, 14).Value = -58837162
$ws.Cells.Item(130, 8).Value = 2676.6667
$ws.Cells.Item(130, 10).Value = 5000
$ws.Cells.Item(130, 12).Value = 15000
$ws.Cells.Item(130, 14).Value = -25040
$ws.Cells.Item(131, 8).Value = 13158777
$ws.Cells.Item(131, 9).Value = 923.63635
$ws.Cells.Item(131, 10).Value = 18519384
$ws.Cells.Item(131, 11).Value = 2770.90905
$ws.Cells.Item(131, 12).Value = 55558152
$ws.Cells.Item(131, 13).Value = 2269.09095
$ws.Cells.Item(131, 14).Value = -55568232
$ws.Cells.Item(136, 8).Value = 1534.1578
$ws.Cells.Item(136, 9).Value = 850
$ws.Cells.Item(136, 10).Value = 3016.5
$ws.Cells.Item(136, 11).Value = 2550
$ws.Cells.Item(136, 12).Value = 9049.5
$ws.Cells.Item(136, 13).Value = 2550
$ws.Cells.Item(136, 14).Value = -19249.5
$ws.Cells.Item(139, 8).Value = 1850.0244
$ws.Cells.Item(139, 9).Value = 1035.4584
$ws.Cells.Item(139, 10).Value = 3000
$ws.Cells.Item(139, 11).Value = 3106.3752
$ws.Cells.Item(139, 12).Value = 9000
$ws.Cells.Item(139, 13).Value = 2033.6248
$ws.Cells.Item(139, 14).Value = -19280
$ws.Cells.Item(140, 8).Value = 4176.8965
$ws.Cells.Item(140, 9).Value = 2375.2632
$ws.Cells.Item(140, 10).Value = 7600
$ws.Cells.Item(140, 11).Value = 7125.7896
$ws.Cells.Item(140, 12).Value = 22800
$ws.Cells.Item(140, 13).Value = -1945.7896
$ws.Cells.Item(140, 14).Value = -33160

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 3115.7222
$ws.Cells.Item(100, 9).Value = 2297.2222
$ws.Cells.Item(100, 10).Value = 3934.2222
$ws.Cells.Item(100, 11).Value = 2297.2222
$ws.Cells.Item(100, 12).Value = 3934.2222
$ws.Cells.Item(100, 13).Value = -1756.2222
$ws.Cells.Item(100, 14).Value = -5016.2222
$ws.Cells.Item(132, 8).Value = 2702.6123
$ws.Cells.Item(132, 9).Value = 3193.9678
$ws.Cells.Item(132, 10).Value = 1856.3889
$ws.Cells.Item(132, 11).Value = 9581.903399999999
$ws.Cells.Item(132, 12).Value = 5569.1667
$ws.Cells.Item(132, 13).Value = -7051.903399999999
$ws.Cells.Item(132, 14).Value = -10629.1667
$ws.Cells.Item(133, 8).Value = 82584.5
$ws.Cells.Item(133, 10).Value = 82584.5
$ws.Cells.Item(133, 12).Value = 82584.5
$ws.Cells.Item(133, 14).Value = -87644.5
$ws.Cells.Item(136, 8).Value = 2296.1562
$ws.Cells.Item(136, 9).Value = 1743.88
$ws.Cells.Item(136, 11).Value = 5231.64
$ws.Cells.Item(136, 13).Value = -2681.64

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 18117008
$ws.Cells.Item(132, 9).Value = 24039506
$ws.Cells.Item(132, 10).Value = 1128.8235
$ws.Cells.Item(132, 11).Value = 72118518
$ws.Cells.Item(132, 12).Value = 3386.4705
$ws.Cells.Item(132, 13).Value = -72115988
$ws.Cells.Item(132, 14).Value = -8446.470499999999
$ws.Cells.Item(136, 8).Value = 949.3421
$ws.Cells.Item(136, 9).Value = 954.4815
$ws.Cells.Item(136, 10).Value = 936.7273
$ws.Cells.Item(136, 11).Value = 2863.4445
$ws.Cells.Item(136, 12).Value = 2810.1819
$ws.Cells.Item(136, 13).Value = -313.4445000000001
$ws.Cells.Item(136, 14).Value = -7910.1819
